# Update the p3_core row (row 8) of the "Modules" manifest sheet for release 1.048.
#
# Columns: A=Module  B=Branch  C=Repo  D=Hash  E=Last Author  F=Last Commit  G=Contributors
#
# p3_core's last commit hash, author, commit timestamp and contributor tally
# all advance to the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHash = "af437528ced4432168cab95af78cd442b5834412"
$newHash = "e6aed994a30ee0dd72784118c164ebe3ff2124a0"

$row = 8

# D8: HYPERLINK formula pointing at the commit - update both the URL and the
# visible label to the new commit hash.
$ws.Cells.Item($row, 4).Formula = '=HYPERLINK("https://github.com/BV-BRC/p3_core/commit/' + $newHash + '", "' + $newHash + '")'

# E8: Last Author
$ws.Cells.Item($row, 5).Value = "olsonanl"

# F8: Last Commit timestamp
$ws.Cells.Item($row, 6).Value = "2025-08-06T20:32:44Z"

# G8: Contributors tally
$ws.Cells.Item($row, 7).Value = "olsonanl:75 bparrello:39 JacobPorter:3"
